$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 498
$ws1.Range("F3").Value = 5889
$ws1.Range("F6").Value = 103
$ws1.Range("F8").Value = 59
$ws1.Range("F10").Value = 27

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 498
$ws4.Range("F3").Value = 5889
$ws4.Range("F7").Value = 103
$ws4.Range("F10").Value = 59
$ws4.Range("F12").Value = 27
